$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added for Cebollín at Terminal La
# Palmera de La Serena. It belongs at the top of the existing data block
# (row 252, right after the header/first block), so insert a fresh row
# there and push the remaining records (old rows 252-341) down by one.
$ws.Rows("252:252").Insert()

$ws.Cells.Item(252, 1).Value = 8
$ws.Cells.Item(252, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(252, 3).Value = "Coquimbo"
$ws.Cells.Item(252, 4).Value = 45119
$ws.Cells.Item(252, 5).Value = 4
$ws.Cells.Item(252, 6).Value = 100112037
$ws.Cells.Item(252, 7).Value = "Cebollín"
$ws.Cells.Item(252, 8).Value = "Sin especificar"
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 1000
$ws.Cells.Item(252, 11).Value = 1000
$ws.Cells.Item(252, 12).Value = 1200
$ws.Cells.Item(252, 13).Value = 1100
$ws.Cells.Item(252, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(252, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(252, 16).Value = 183
$ws.Cells.Item(252, 17).Value = 6
$ws.Cells.Item(252, 18).Value = "Hortaliza"
